$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "CSS Id" column header to "VLJ #"
$ws.Range("C2").Value = "VLJ #"

# Judge Bernard Jones' id changes from the old CSS login id to a VLJ number
$ws.Range("C3").Value = "123"
$ws.Range("C4").Value = "123"
$ws.Range("C5").Value = "123"
$ws.Range("C6").Value = "123"
$ws.Range("C7").Value = "123"

# Second example judge is swapped out: Roth, Lauren (DSUSER) -> Huels, Stuart (860)
$ws.Range("B8").Value = "Huels, Stuart"
$ws.Range("B9").Value = "Huels, Stuart"
$ws.Range("C8").Value = "860"
$ws.Range("C9").Value = "860"

# Extend the table with one additional blank data row (row 10), matching the
# look (borders/fill/height) of the last existing row.
$ws.Range("A9:I9").Copy()
$ws.Range("A10:I10").PasteSpecial(-4122)
$ws.Rows(10).RowHeight = 17
$ws.Range("A10:I10").ClearContents()
